$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

# Rename sheet "Part Codes" -> "source"
$ws.Name = "source"

# Update descriptions (column C) for specific PN rows
$ws.Range("C63").Value = "Z Motion, L Limit Switch Mount"
$ws.Range("C64").Value = "Z Motion, R Limit Switch Mount"
$ws.Range("C91").Value = "Electrical, Buck Mount, DROK 5A"
$ws.Range("C111").Value = "Skirt, Hexagon Grill"

# Update the view: change the active selection (also resets the scrolled
# top-left cell back to the sheet's default origin)
$ws.Range("I156").Select()
